$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.198.28"
$ws.Range("E2").Value = "  -0.06%  "
$ws.Range("D3").Value = "1.655.01"
$ws.Range("E3").Value = "  -0.31%  "
$ws.Range("E4").Value = "  +0.45%  "
$ws.Range("D5").Value = "218.58"
$ws.Range("E5").Value = "  +1.38%  "
$ws.Range("D6").Value = "0.5220"
$ws.Range("E6").Value = "  -0.75%  "
$ws.Range("E7").Value = "  +0.41%  "
$ws.Range("D8").Value = "0.2652"
$ws.Range("E8").Value = "  +0.65%  "
$ws.Range("D9").Value = "0.06317"
$ws.Range("E9").Value = "  -1.51%  "
$ws.Range("E10").Value = "  +0.68%  "
$ws.Range("D11").Value = "0.07744"
$ws.Range("E11").Value = "  -0.51%  "
$ws.Range("D12").Value = "1.654.06"
$ws.Range("E12").Value = "  -0.49%  "
$ws.Range("D13").Value = "4.423"
$ws.Range("E13").Value = "  -1.03%  "
$ws.Range("D14").Value = "0.5463"
$ws.Range("E14").Value = "  -1.49%  "
$ws.Range("D15").Value = "0.0₅8210"
$ws.Range("E15").Value = "  -1.04%  "
$ws.Range("D16").Value = "64.85"
$ws.Range("E16").Value = "  -0.77%  "
$ws.Range("D17").Value = "26.207.94"
$ws.Range("E17").Value = "  -0.01%  "
$ws.Range("D18").Value = "1.006"
$ws.Range("E18").Value = "  +0.42%  "
$ws.Range("E19").Value = "  -2.04%  "
$ws.Range("D20").Value = "192.55"
$ws.Range("E20").Value = "  +0.48%  "
$ws.Range("D21").Value = "10.17"
$ws.Range("E21").Value = "  -1.19%  "
$ws.Range("D22").Value = "6.140"
$ws.Range("E22").Value = "  -3.43%  "
$ws.Range("D23").Value = "1.007"
$ws.Range("E23").Value = "  +0.59%  "
$ws.Range("D24").Value = "138.71"
$ws.Range("E24").Value = "  -2.93%  "
$ws.Range("D25").Value = "0.1240"
$ws.Range("E25").Value = "  -1.81%  "
$ws.Range("D26").Value = "7.264"
$ws.Range("E26").Value = "  -2.36%  "
$ws.Range("D27").Value = "16.14"
$ws.Range("E27").Value = "  +0.49%  "
$ws.Range("D28").Value = "1.416"
$ws.Range("E28").Value = "  -1.11%  "
$ws.Range("D29").Value = "0.06078"
$ws.Range("E29").Value = "  -0.19%  "
$ws.Range("D30").Value = "1.285"
$ws.Range("E30").Value = "  +1.35%  "
$ws.Range("D31").Value = "3.549"
$ws.Range("E31").Value = "  -0.25%  "
$ws.Range("D32").Value = "3.359"
$ws.Range("E32").Value = "  -2.20%  "
$ws.Range("D33").Value = "1.655"
$ws.Range("E33").Value = "  -0.52%  "
$ws.Range("D34").Value = "0.9856"
$ws.Range("E34").Value = "  -1.80%  "
$ws.Range("E35").Value = "  +0.42%  "
$ws.Range("D36").Value = "2.771"
$ws.Range("E36").Value = "  +0.16%  "
$ws.Range("D37").Value = "0.5972"
$ws.Range("E37").Value = "  +4.79%  "
$ws.Range("D38").Value = "0.01598"
$ws.Range("E38").Value = "  -0.63%  "
$ws.Range("D39").Value = "5.972"
$ws.Range("E39").Value = "  +0.81%  "
$ws.Range("D40").Value = "0.8606"
$ws.Range("E40").Value = "  +0.31%  "
$ws.Range("D41").Value = "1.056.33"
$ws.Range("E41").Value = "  +2.17%  "
$ws.Range("D42").Value = "1.003"
$ws.Range("E42").Value = "  +0.20%  "
$ws.Range("D43").Value = "99.86"
$ws.Range("E43").Value = "  +0.09%  "
$ws.Range("D44").Value = "1.793.09"
$ws.Range("E44").Value = "  -0.76%  "
$ws.Range("D45").Value = "0.0₈109"
$ws.Range("E45").Value = "  -0.31%  "
$ws.Range("D46").Value = "57.10"
$ws.Range("E46").Value = "  +1.31%  "
$ws.Range("D47").Value = "1.006"
$ws.Range("E47").Value = "  +0.42%  "
$ws.Range("D48").Value = "8.041"
$ws.Range("E48").Value = "  -0.59%  "
$ws.Range("D49").Value = "0.05181"
$ws.Range("E49").Value = "  +0.25%  "
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").Value = "1.477"
$ws.Range("E50").Value = "  +4.86%  "
$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").Value = "0.4231"
$ws.Range("E51").Value = "  +0.35%  "
